$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 231, pushing existing rows 231-335 down to 232-336
$ws.Rows.Item(231).Insert()

# Populate the newly inserted row 231 with its data
$ws.Range("A231").Value = 5
$ws.Range("B231").Value = "Macroferia Regional de Talca"
$ws.Range("C231").Value = "Maule"
$ws.Range("D231").Value = 44609
$ws.Range("E231").Value = 7
$ws.Range("F231").Value = 100112043
$ws.Range("G231").Value = "Pepino ensalada"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 300
$ws.Range("K231").Value = 12000
$ws.Range("L231").Value = 12000
$ws.Range("M231").Value = 12000
$ws.Range("N231").Value = "$/caja 80 unidades"
$ws.Range("O231").Value = "Región del Maule"
$ws.Range("P231").Value = 150
$ws.Range("Q231").Value = 80
$ws.Range("R231").Value = "Hortaliza"

# Apply the same date style (numFmt) as the date column in surrounding rows
$ws.Range("D231").NumberFormat = $ws.Range("D232").NumberFormat
